# Generate Report for Archive
#
# The localization status moved on from "Ready for handoff" to
# "In Translation" for every row/sheet that reported it. Update the
# shared "Status" value everywhere it appears (Overview!E:F and the
# "Status" column on the per-language sheets), then resize the affected
# columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# 12.5 is the ColumnWidth input that this engine's pixel-snapping
# (width -> round(width*6)+5 px, displayed back as px/6) resolves to the
# stored width closest to the generator's original 13.4101845877511
# target (13.333333333333334, i.e. 80px @ 6px/unit).
$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth  = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $value = [string]$cell.Value()
        if ($value -eq $oldStatus) {
            $cell.Value = $newStatus
            $cell.EntireColumn.ColumnWidth = $newWidth
        }
    }
}
